$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")
$wsDesc = $wb.Worksheets.Item("Description")

# Target data for rows 2-65 (A,B,C,D,E,F,G,H,I,J,K,L), comma-separated.
# This reflects the "Updated LRLC Runs and Results" commit:
#  - Electric Total columns (B=kWh_Elec_Total, D=kTDV_Total_Elec_2022, F=Source_kBtu_Elec)
#    carry freshly simulated values (PVB Removed from the electric total calc).
#  - All other metric columns (C,E,G Total-NG + H,I,J,K Comp) and the Measure label (L)
#    are reassigned as a block per CZ row to the new Measure grouping.
$data = @"
1,22412.338059469,2109.97350786,88.2884163202638,102.510602155379,59114.4892718639,184792.390822377,2529.16017914569,1816.07737098014,9.51917065356799,88.3860451304422,lrgs-low-baseszac
2,22543.519555299,1580.90003854,92.6769216449942,77.3321558462802,59800.5719615755,138455.90794611,2521.11881163574,1287.0038841677,15.6279803108575,63.2075979046658,lrgs-low-baseszac
3,22162.690861452,1374.46742384,88.60930678235,66.72803064856,58644.3115179705,120376.450421159,2171.42629957606,1080.5712609843,10.5013610078282,52.6034721607589,lrgs-low-baseszac
4,23103.489031248,1281.7087638,96.2718150811335,62.3227211764989,61743.3095376734,112252.606925296,2965.21171576757,987.81262856462,19.1127050617248,48.198164044124,lrgs-low-baseszac
5,22153.317999506,1330.1573014,87.49378582589,64.3889010001654,58582.3834437173,116495.750766486,2142.06295214739,1036.26116291667,8.60612409327883,50.2643437743682,lrgs-low-baseszac
6,22824.736391637,1094.31645402,94.8985928965348,52.8628290105597,61068.3480344425,94492.9330530789,2625.31581731824,800.4202732668,16.8807864172737,38.7155745127215,lrgs-low-baseszac
7,22972.210245018,1089.44321467,93.2347998011782,52.6458675040139,61667.3231348664,94072.1345921222,2746.2528246686,795.54703899072,14.5616778645056,38.4647188794018,lrgs-low-baseszac
8,24891.113121748,1041.96822101,103.665282581345,50.3130887356167,67358.1377059218,89972.7249733323,4581.23745280582,748.07204671984,27.6021779609433,36.1658346612252,lrgs-low-baseszac
9,24637.782387486,1075.12231205,103.402841217832,52.0264230709302,66766.1146227962,92835.5415686324,4368.3335429612,781.2261764717,27.479545599953,37.8791709070342,lrgs-low-baseszac
10,25699.259990046,1120.20199518,110.578477409453,54.456515719586,69676.7140988046,96728.1189528149,5426.9141796189,826.3058358483,34.0138954295134,40.309262321822,lrgs-low-baseszac
11,27196.29254902,1421.09476549,122.230661607986,69.837867050281,73012.0443149839,124460.093134727,6988.90714459473,1127.19865573481,44.6179577092243,55.7133112090727,lrgs-low-baseszac
12,24091.675496012,1433.67323173,106.063459218684,70.2604368080786,64778.6269724753,125561.720638916,3942.46287148064,1139.77708303009,28.5464985966722,56.1358790510155,lrgs-low-baseszac
13,27895.3627617635,1293.96971283,123.036424866407,63.4464474354856,75546.1545320567,113326.426135141,7645.8934712822,1000.07356264392,45.5448177633257,49.3218895840806,lrgs-low-baseszac
14,26492.45677319,1359.74079552,110.557331935402,66.834128246961,71822.0066868266,117412.011387214,6296.27529205534,1065.84462763463,35.2806189234038,52.6868745376419,lrgs-low-baseszac
15,35387.785734647,847.71650831,151.824172414797,41.1137036774873,93769.1497170609,73199.3190575408,14925.0356490264,553.8203714833,75.1518529256129,26.9664513892842,lrgs-low-baseszac
16,23297.692842667,2189.08371792,93.275593705566,107.632350620382,62385.293211311,189024.793006739,3253.3613152666,1895.1876262228,13.3046884508655,93.4851007054751,lrgs-low-baseszac
1,149753.982479712,8696.770195,109.294308175781,78.2650760720572,387593.56732525,761666.888603171,42586.8161831916,7352.46149568424,31.8627964983222,66.2504924891521,lrlc-baseline-all
2,143041.6992896,7341.997463,108.997595565136,66.2800542546027,373325.813865395,643015.307796756,35232.7420231726,5997.6890156942,32.916111155739,54.2654728395149,lrlc-baseline-all
3,136229.96591973,6625.982674,100.565010403691,59.4585968502382,353206.967151864,580306.423428968,28498.1004077125,5281.6741805551,23.5884503389445,47.4440151078061,lrlc-baseline-all
4,143535.24887464,6375.030487,111.307895238578,57.279707155273,373253.095152146,558327.922540173,35134.3757147033,5030.7218267762,35.2501284085163,45.265123941294,lrlc-baseline-all
5,136820.48302102,6554.546615,99.7372002783615,58.7383453960818,355026.99625769,574050.022538453,28949.4752215308,5210.237931358,22.1835100841028,46.7237619116295,lrlc-baseline-all
6,139719.99736543,5810.479147,107.586362881489,52.0435003408487,360841.743775931,501728.010235828,30949.3851436631,4466.1702918795,30.6945074197909,40.0095705356774,lrlc-baseline-all
7,137992.0018959,5760.495689,103.36197207864,51.6950752895751,357931.598064233,497411.99768461,29093.0850377257,4416.187017021,25.9811005159642,39.6322285929515,lrlc-baseline-all
8,149113.68023369,5649.62038,114.130295366572,50.614321783342,383473.062026314,487838.045733061,39993.7012375079,4305.311699464,38.9201490411927,38.5803935371693,lrlc-baseline-all
9,149508.53552941,5770.703551,115.628263460512,51.7548839427498,386228.761938985,498293.434509431,40566.4232453719,4426.394947388,40.4635656667657,39.7209563184306,lrlc-baseline-all
10,157035.044007763,5928.14455,122.696684221462,53.3281616110487,405896.3527276,511888.278782919,48118.6978464883,4583.835910441,47.3143877531963,41.2942336632156,lrlc-baseline-all
11,170484.51045498,6840.718989,137.275071783378,61.9669056801819,442736.873909379,599113.122611406,61895.6682311594,5496.4104890889,60.9481537796566,49.95232388295,lrlc-baseline-all
12,152637.976425711,6914.970738,122.523691735952,62.4991409877888,399849.078186213,605616.122847826,44251.9394230849,5570.662225892,46.290917311368,50.4845590993907,lrlc-baseline-all
13,173851.647550698,6393.784235,139.759897715118,57.7931467903642,451482.82266961,559970.383887148,65045.247717797,5049.4757180326,63.4976705978154,45.7785648282996,lrlc-baseline-all
14,167366.62840467,6754.363389,128.18728039526,61.2338668307596,439591.284920256,583231.299491435,58819.8813587715,5410.0548007928,53.792049413419,49.1999394187638,lrlc-baseline-all
15,209942.47879774,4919.170673,164.413306323774,44.2103882820078,516032.561230557,424764.57643756,100229.835889598,3574.8619825,89.0403794961183,32.1764598235932,lrlc-baseline-all
16,158004.16494594,9471.761148,117.849214719056,86.0533486654822,414267.954531748,817875.385830907,50065.887172835,8127.4525328601,39.5226900278527,74.0194207932843,lrlc-baseline-all
1,30735.172454769,1328.78917726,122.778491963397,63.9254211578953,86083.9336122595,116375.929863603,10825.3347131744,1034.89303802,43.9068809473274,49.8008639104389,lrgs-low-base
2,27330.246599304,1181.25209414,111.747470904478,56.9669867387733,77463.5536861589,103454.568423214,7293.473423777,887.35592234,34.6422812339047,42.8424278851332,lrgs-low-base
3,24385.710717976,1180.28488744,97.808730748672,56.8230823603711,67641.3929961064,103369.860042826,4384.8393753676,886.38872916,19.6621832883062,42.6985241518383,lrgs-low-base
4,25294.326779408,1107.1874831,107.347662692801,53.4144432284655,70405.4090751455,96967.9578101296,5148.134068229,813.29133789,30.1613083051692,39.289885624116,lrgs-low-base
5,24464.15438881,1161.1911062,97.1909158330939,55.8479862847407,67653.3689246836,101697.618437879,4441.1625472791,867.2949263,18.2583930264168,41.7234270269359,lrgs-low-base
6,23394.621338671,1058.35661922,97.6882410478758,51.0161440534564,63443.3536028196,91387.843798619,3191.4987101027,764.46045407,19.6558107601317,36.8688903672142,lrgs-low-base
7,23234.731644616,1075.48238307,94.5869859119648,51.9219771661029,62782.1994446294,92866.6332758458,3007.114908791,781.58620879,15.9073743682145,37.7408286144227,lrgs-low-base
8,25400.320883739,1026.47003241,106.338195006113,49.5157372603914,69369.3520135814,88634.4746962355,5088.2789614688,732.57386334,30.26670334708,35.3684834509751,lrgs-low-base
9,25552.200325904,1026.00030815,107.720522500849,49.5092849118198,70415.5100839497,88593.914561285,5278.3648292043,732.10416345,31.779855167847,35.3620322682811,lrgs-low-base
10,27401.114402095,1019.49318628,118.288615945215,49.2692583060921,76646.2760011499,88032.0322748848,7122.1085625549,725.59702795,41.6979987959261,35.1220049507782,lrgs-low-base
11,31774.69057071,1069.76081079,141.387224394219,51.8209032565899,90866.036016694,93690.1136898472,11555.9107257463,775.86468003,63.7279721857703,37.6963463568015,lrgs-low-base
12,28079.633851031,1107.70706023,123.478478549197,53.5479610122171,80574.649564886,97013.4625995081,7919.531663948,813.81091119,45.9175732204687,39.4234032395588,lrgs-low-base
13,31437.9920177735,1034.58997803,138.271901604602,50.1147803875143,89535.5868956426,90609.836971337,11179.6821328472,740.69383555,60.7424368239093,35.9902229052603,lrgs-low-base
14,31349.14673276,1030.46181732,132.398443399189,49.920815174173,91465.5303112534,88979.1606075889,11140.7236489092,736.56565708,57.0779875603305,35.7735617841142,lrgs-low-base
15,36348.027145147,843.41277021,156.534115111584,40.8915209502375,96683.86451116,72827.6963567513,15884.4790040952,549.51663448,79.8586770822312,26.744268718725,lrgs-low-base
16,34983.889619326,1247.39889322,148.536486762246,60.3360520595744,99583.0427430117,107711.420836745,14916.3700450375,953.50272313,68.4760140985721,46.1887981769321,lrgs-low-base
1,174539.349900922,6619.330973,128.510969462808,59.1067270304137,469390.159800459,579723.864583444,67286.8278348038,5275.022421,51.0179576622851,47.0921446191635,lrlc-flat19-base
2,158415.57558209,6186.758084,120.557879089585,55.31458802475,430516.548583787,541838.944196776,50560.2932897283,4842.449626,44.4427335927938,43.3000065195435,lrlc-flat19-base
3,142451.48789468,6176.655733,105.472106669281,55.1829882749937,378741.555426775,540954.175934396,34692.2513553354,4832.347197,28.4750778340043,43.1684061377285,lrlc-flat19-base
4,150489.67179329,5916.584654,117.820173635166,52.9186196440184,400774.212619399,518177.038547061,42064.6325372254,4572.276035,41.746400837232,40.90403674999,lrlc-flat19-base
5,143294.10220174,6177.467889,104.914213405078,55.1824275327287,380816.650967763,541025.304907534,35390.0035800967,4833.159286,27.3372585873033,43.1678448063354,lrlc-flat19-base
6,141464.11055318,5760.430606,109.372463989217,51.5664691108899,367766.217721587,497406.377844445,32684.9360938539,4416.12169,32.474559590281,39.5325387384548,lrlc-flat19-base
7,138942.52003369,5747.447556,104.422927870384,51.5693343359411,361725.655180016,496285.306814243,30041.2095182028,4403.138861,27.040629985282,39.5064874174496,lrlc-flat19-base
8,151112.69943237,5624.593541,116.114020802979,50.3746175429341,390983.257279837,485677.007750429,41987.6102970714,4280.284804,40.9002613331147,38.3406887624799,lrlc-flat19-base
9,152666.5324993,5663.812742,118.499188828883,50.7373254036727,398442.398152766,489063.539425863,43712.7240025002,4319.504178,43.3261238396236,38.7033981521308,lrlc-flat19-base
10,163074.344820418,5647.063745,127.867230108075,50.6348891998656,430199.897034407,487617.283321047,54135.973846145,4302.755109,52.4689615473748,38.6009612904814,lrlc-flat19-base
11,186418.26606348,5744.446926,149.771795109047,51.4989934475539,504709.25768682,503101.142006486,77791.0292286628,4400.138424,73.4157137815196,39.4844116325389,lrlc-flat19-base
12,166329.139744071,5917.944872,134.026400391812,52.9795177173176,453843.06440615,518296.16702679,57907.2156753563,4573.636293,57.7671898730524,40.9649352029195,lrlc-flat19-base
13,186159.330271138,5619.207819,149.706149466242,50.3839876701363,499313.442393003,492132.64694208,77323.7347515964,4274.89929,73.4206159768654,38.3694056113146,lrlc-flat19-base
14,183676.64595247,5773.985992,142.149458381102,51.860282639682,505310.158402147,498576.869412127,75088.843026052,4429.677375,67.7271809432991,39.8263549331731,lrlc-flat19-base
15,214149.49953893,4916.777897,168.230102446458,44.1874367512114,528144.410534121,424557.963056624,104436.356307479,3572.469213,92.8568541236453,32.1535083549023,lrlc-flat19-base
16,197096.38234033,6529.732686,152.766732696951,58.5492934473676,538667.296211295,563834.703650926,89082.5830168385,5185.424042,74.3867678145332,46.5153654105036,lrlc-flat19-base
"@

$lines = $data -split "`r?`n" | Where-Object { $_.Trim().Length -gt 0 }

$row = 2
foreach ($line in $lines) {
    $p = $line.Trim() -split ","

    $ws.Cells.Item($row, 1).Value  = [double]$p[0]   # A: CZ
    $ws.Cells.Item($row, 2).Value  = [double]$p[1]   # B: kWh_Elec_Total
    $ws.Cells.Item($row, 3).Value  = [double]$p[2]   # C: Therm_NG_Total
    $ws.Cells.Item($row, 4).Value  = [double]$p[3]   # D: kTDV_Total_Elec_2022
    $ws.Cells.Item($row, 5).Value  = [double]$p[4]   # E: kTDV_Total_NG_2022
    $ws.Cells.Item($row, 6).Value  = [double]$p[5]   # F: Source_kBtu_Elec
    $ws.Cells.Item($row, 7).Value  = [double]$p[6]   # G: Source_kBtu_NG
    $ws.Cells.Item($row, 8).Value  = [double]$p[7]   # H: kWh_Elec_Comp
    $ws.Cells.Item($row, 9).Value  = [double]$p[8]   # I: Therm_NG_Comp
    $ws.Cells.Item($row, 10).Value = [double]$p[9]   # J: kTDV_Comp_Elec_2022
    $ws.Cells.Item($row, 11).Value = [double]$p[10]  # K: kTDV_Comp_NG_2022
    $ws.Cells.Item($row, 12).Value = $p[11]          # L: Measure

    $row = $row + 1
}

# The shared "description" string (Description sheet, A1) gains a trailing
# "PVB Removed: ,TRUE" clause.
$wsDesc.Range("A1").Value = "LifeTime: ,30,TDVMultiplier: ,res,DHW Removed: ,FALSE,PVB Removed: ,TRUE"

Write-Output "rows written: $($row - 2)"
